$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.985
$ws.Range("C6").Value = -11.626
$ws.Range("A14").Value = -21.852
$ws.Range("B15").Value = 5.872
$ws.Range("A16").Value = -22.225
$ws.Range("C18").Value = -12.237
$ws.Range("C19").Value = -12.018
$ws.Range("A21").Value = -21.571
$ws.Range("B21").Value = 6.891000000000001
$ws.Range("B22").Value = 7.358
$ws.Range("A23").Value = -20.916
$ws.Range("B24").Value = 5.378
$ws.Range("A25").Value = -21.683
$ws.Range("A26").Value = -21.398
$ws.Range("B27").Value = 6.299000000000001
$ws.Range("B28").Value = 5.827
$ws.Range("A29").Value = -21.498
$ws.Range("C35").Value = -12.446
$ws.Range("B36").Value = 7.150999999999999
$ws.Range("B39").Value = 7.311
$ws.Range("A40").Value = -20.624
$ws.Range("C44").Value = -12.275
$ws.Range("B45").Value = 5.451000000000001
$ws.Range("C47").Value = -12.454
$ws.Range("B48").Value = 5.48
$ws.Range("B49").Value = 6.470000000000001
$ws.Range("C50").Value = -13.233
$ws.Range("C51").Value = -11.09
$ws.Range("B52").Value = 5.272
$ws.Range("C52").Value = -10.877
$ws.Range("A53").Value = -21.122
$ws.Range("B53").Value = 6.327
$ws.Range("B54").Value = 5.666
$ws.Range("C55").Value = -13.143
$ws.Range("A57").Value = -21.619
$ws.Range("B57").Value = 6.093
$ws.Range("C57").Value = -13.361
$ws.Range("C58").Value = -12.59
$ws.Range("A59").Value = -22.574
$ws.Range("C64").Value = -10.885
$ws.Range("A65").Value = -21.515
$ws.Range("C66").Value = -11.507
$ws.Range("A69").Value = -21.614
$ws.Range("B70").Value = 4.760000000000001
$ws.Range("B71").Value = 4.972
$ws.Range("A79").Value = -21.195
$ws.Range("C80").Value = -12.335
$ws.Range("A83").Value = -21.778
$ws.Range("C83").Value = -12.391
$ws.Range("B86").Value = 5.613000000000001
$ws.Range("B87").Value = 5.766000000000001
$ws.Range("B89").Value = 5.747999999999999
$ws.Range("A91").Value = -21.111
$ws.Range("C92").Value = -10.888
$ws.Range("A93").Value = -21.687
$ws.Range("C94").Value = -10.924
$ws.Range("C96").Value = -10.867
$ws.Range("C97").Value = -11.424
$ws.Range("A100").Value = -22.615
$ws.Range("B101").Value = 6.084000000000001
$ws.Range("C101").Value = -12.879
$ws.Range("A103").Value = -21.894
